# Applies the "Updated cryptos list" edit: refreshed Price/Volume(1h)
# figures for most rows, plus a name/link/price swap between the
# dogwifhat and FirstDigitalUSD rows (45 <-> 46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the write to land as text (matches the source data, which
    # stores numbers like "557.84" or "1.00" as literal strings) -- a
    # plain .Value assignment would let Excel auto-coerce numeric-looking
    # strings into real numbers (dropping formatting like trailing zeros).
    # Restoring the style to "Normal" afterwards avoids leaving a stray
    # text-format style on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.470.03"
Set-TextValue $ws.Range("E2") "  +1.66%  "
Set-TextValue $ws.Range("D3") "3.587.30"
Set-TextValue $ws.Range("E3") "  +0.52%  "
Set-TextValue $ws.Range("E4") "  +0.20%  "
Set-TextValue $ws.Range("D5") "200.13"
Set-TextValue $ws.Range("E5") "  +8.01%  "
Set-TextValue $ws.Range("D6") "557.84"
Set-TextValue $ws.Range("E6") "  -4.95%  "
Set-TextValue $ws.Range("D7") "3.582.52"
Set-TextValue $ws.Range("E7") "  +0.63%  "
Set-TextValue $ws.Range("E8") "  +0.01%  "
Set-TextValue $ws.Range("E9") "  +0.16%  "
Set-TextValue $ws.Range("D10") "0.670"
Set-TextValue $ws.Range("E10") "  +0.09%  "
Set-TextValue $ws.Range("D11") "59.33"
Set-TextValue $ws.Range("E11") "  +10.97%  "
Set-TextValue $ws.Range("D12") "0.151"
Set-TextValue $ws.Range("E12") "  +3.43%  "
Set-TextValue $ws.Range("D13") "0.0000286"
Set-TextValue $ws.Range("E13") "  +11.72%  "
Set-TextValue $ws.Range("D14") "9.96"
Set-TextValue $ws.Range("E14") "  +2.04%  "
Set-TextValue $ws.Range("D15") "4.175.45"
Set-TextValue $ws.Range("E15") "  +0.82%  "
Set-TextValue $ws.Range("D16") "3.594.70"
Set-TextValue $ws.Range("E16") "  +0.84%  "
Set-TextValue $ws.Range("D18") "18.94"
Set-TextValue $ws.Range("E18") "  +3.73%  "
Set-TextValue $ws.Range("D19") "67.412.12"
Set-TextValue $ws.Range("E19") "  +1.79%  "
Set-TextValue $ws.Range("E21") "  +1.68%  "
Set-TextValue $ws.Range("D22") "398.64"
Set-TextValue $ws.Range("E22") "  +0.69%  "
Set-TextValue $ws.Range("D23") "12.81"
Set-TextValue $ws.Range("E23") "  +15.39%  "
Set-TextValue $ws.Range("E24") "  -5.57%  "
Set-TextValue $ws.Range("D25") "84.92"
Set-TextValue $ws.Range("E25") "  -0.55%  "
Set-TextValue $ws.Range("D26") "2.93"
Set-TextValue $ws.Range("E26") "  +2.18%  "
Set-TextValue $ws.Range("E27") "  +0.33%  "
Set-TextValue $ws.Range("D28") "3.86"
Set-TextValue $ws.Range("E28") "  +8.77%  "
Set-TextValue $ws.Range("D29") "6.10"
Set-TextValue $ws.Range("E29") "  +1.11%  "
Set-TextValue $ws.Range("D30") "8.36"
Set-TextValue $ws.Range("E30") "  +19.28%  "
Set-TextValue $ws.Range("D31") "9.44"
Set-TextValue $ws.Range("E31") "  +5.42%  "
Set-TextValue $ws.Range("D32") "31.37"
Set-TextValue $ws.Range("E32") "  +1.44%  "
Set-TextValue $ws.Range("D33") "663.14"
Set-TextValue $ws.Range("E33") "  +7.37%  "
Set-TextValue $ws.Range("E34") "  +0.28%  "
Set-TextValue $ws.Range("D35") "63.68"
Set-TextValue $ws.Range("E35") "  +0.48%  "
Set-TextValue $ws.Range("E36") "  +0.52%  "
Set-TextValue $ws.Range("D37") "42.21"
Set-TextValue $ws.Range("D38") "0.427"
Set-TextValue $ws.Range("E38") "  +8.46%  "
Set-TextValue $ws.Range("D39") "1.00"
Set-TextValue $ws.Range("D40") "0.0₃0767"
Set-TextValue $ws.Range("E40") "  +2.04%  "
Set-TextValue $ws.Range("D41") "3.19"
Set-TextValue $ws.Range("E41") "  +13.73%  "
Set-TextValue $ws.Range("D42") "3.253.29"
Set-TextValue $ws.Range("E42") "  +8.16%  "
Set-TextValue $ws.Range("E43") "  +3.60%  "
Set-TextValue $ws.Range("D44") "2.78"
Set-TextValue $ws.Range("E44") "  +12.09%  "
Set-TextValue $ws.Range("B45") "FirstDigitalUSD"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D45") "0.999"
Set-TextValue $ws.Range("E45") "  -0.05%  "
Set-TextValue $ws.Range("B46") "dogwifhat"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D46") "3.00"
Set-TextValue $ws.Range("E46") "  +28.75%  "
Set-TextValue $ws.Range("D47") "0.0415"
Set-TextValue $ws.Range("E47") "  +2.34%  "
Set-TextValue $ws.Range("D48") "2.75"
Set-TextValue $ws.Range("E48") "  +10.83%  "
Set-TextValue $ws.Range("D49") "3.11"
Set-TextValue $ws.Range("E49") "  +2.19%  "
Set-TextValue $ws.Range("E50") "  +0.20%  "
Set-TextValue $ws.Range("D51") "8.70"
Set-TextValue $ws.Range("E51") "  +1.29%  "

Write-Host "Applied cryptos list update ($($wb.Name))"
